$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells (L1:N1) - copy style from K1 (existing header style)
$ws.Range("L1").Value = "apoio_medio"
$ws.Range("M1").Value = "contribuicoes"
$ws.Range("N1").Value = "media_contribuicoes"

$ws.Range("K1").Copy()
$ws.Range("L1:N1").PasteSpecial(-4122)

# New data values for rows 2-6, columns L (apoio_medio), M (contribuicoes), N (media_contribuicoes)
$ws.Range("L2").Value = 91.30338044842225
$ws.Range("M2").Value = 263553
$ws.Range("N2").Value = 317.533734939759

$ws.Range("L3").Value = 0
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0

$ws.Range("L4").Value = 90.1669167946294
$ws.Range("M4").Value = 203646
$ws.Range("N4").Value = 147.2494577006508

$ws.Range("L5").Value = 19.17132323902399
$ws.Range("M5").Value = 2063
$ws.Range("N5").Value = 15.05839416058394

$ws.Range("L6").Value = 25.0794338805401
$ws.Range("M6").Value = 145
$ws.Range("N6").Value = 9.666666666666666
